$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values first (B1 = 0, A2 = 0, B2 = text).
$ws.Range("B1").Value = 0
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "disconnected_elements"

# Build the shared "bold, centered/top, thin-boxed" look once on a scratch
# cell so the style table doesn't pick up orphaned intermediate xfs, then
# stamp it onto B1 and A2 via a formats-only paste.
$scratch = $ws.Range("Z1")
$scratch.Font.Bold = $true
$scratch.Borders.LineStyle = 1       # xlContinuous
$scratch.Borders.Weight = 2          # xlThin
$scratch.HorizontalAlignment = -4108 # xlCenter
$scratch.VerticalAlignment = -4160   # xlTop

$scratch.Copy()
$ws.Range("B1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A2").PasteSpecial(-4122)  # xlPasteFormats

$scratch.Clear()
$excel.CutCopyMode = $false
